$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("B11").Value = "Potato Blight disease observation"

# Row 12
$ws.Range("B12").Value = "https://orcid.org/0000-0002-6665-0896"
$ws.Range("C12").Value = "Poul Lassen"

# Row 13
$ws.Range("B13").Value = "https://orcid.org/0000-0001-7682-5613"
$ws.Range("C13").Value = "Jens Grønbech Hansen"

# Row 14
$ws.Range("B14").Value = "https://orcid.org/0000-0002-6323-4942"
$ws.Range("C14").Value = "Margit Styrbæk Jørgensen"

# Row 15
$ws.Range("B15").Value = "https://orcid.org/0000-0003-4093-2147"
$ws.Range("C15").Value = "Eva Overby Bach"

# Row 16
$ws.Range("A16").Value = "dct:creator"
$ws.Range("B16").Value = "https://orcid.org/0000-0002-0721-551X"
$ws.Range("C16").Value = "Ying Wang"

# Row 17
$ws.Range("A17").Value = "dct:creator"
$ws.Range("C17").Value = ""

# Row 18
$ws.Range("A18").Value = "dct:rights"
$ws.Range("B18").Value = "https://spdx.org/licenses/CC-BY-4.0.html"
$ws.Range("C18").Value = "License under which the vocabulary is provided"

# Row 19
$ws.Range("A19").Value = "pav:version"
$ws.Range("B19").Value = "0.0.1"
$ws.Range("C19").Value = "Vocabulary version"

# Row 20
$ws.Range("A20").Value = "pav:createdOn"
$ws.Range("B20").Value = "2022-01-21T10:03:28Z"
$ws.Range("C20").Value = "Date when vocabulary was initially created (follow https://en.wikipedia.org/wiki/ISO_8601)"

# Row 21
$ws.Range("A21").Value = "pav:lastUpdatedOn"
$ws.Range("B21").Value = "2022-01-21T10:03:28Z"
$ws.Range("C21").Value = "Date of the last vocabulary update"
$ws.Range("D21").Value = ""
$ws.Range("E21").Value = ""
$ws.Range("F21").Value = ""
$ws.Range("G21").Value = ""
$ws.Range("H21").Value = ""
$ws.Range("I21").Value = ""
$ws.Range("J21").Value = ""
$ws.Range("K21").Value = ""
$ws.Range("L21").Value = ""

# Row 22
$ws.Range("A22").Value = "Definition of terms (optionally properties)"
$ws.Range("B22").Value = ""
$ws.Range("E22").Value = ""

# Row 23
$ws.Range("A23").Value = "Identifier"
$ws.Range("B23").Value = "skos:prefLabel@en"
$ws.Range("C23").Value = 'qudt:unit(separator=",")'
$ws.Range("D23").Value = 'skos:altLabel(separator=";")'
$ws.Range("E23").Value = "skos:definition@en"
$ws.Range("F23").Value = 'dct:source(separator=",")'
$ws.Range("G23").Value = 'skos:broader(separator=",")'
$ws.Range("H23").Value = 'skos:exactMatch(separator=",")'
$ws.Range("I23").Value = 'skos:closeMatch(separator=",")'
$ws.Range("J23").Value = "skos:editorialNote@en"
$ws.Range("K23").Value = 'dct:creator(separator=",")'
$ws.Range("L23").Value = 'dct:contributor(separator=",")'

# Row 24
$ws.Range("A24").Value = "vars:SampleID"
$ws.Range("B24").Value = "SampleID"
$ws.Range("E24").Value = "Number given to the sample after institutes nameing standard"

# Row 25
$ws.Range("A25").Value = "vars:ObservationID"
$ws.Range("B25").Value = "ObservationID"
$ws.Range("E25").Value = "Running number"

# Row 26
$ws.Range("A26").Value = "vars:CropSeasonYear"
$ws.Range("B26").Value = "CropSeasonYear"
$ws.Range("E26").Value = "4 digit number representing the year the disease was observed. Automatically created from the Blight Tracker App."

# Row 27
$ws.Range("A27").Value = "vars:CountryCode"
$ws.Range("B27").Value = "CountryCode"
$ws.Range("E27").Value = "2 char code as defined in the ISO standard ISO 3166-1 alpha-2 codes "
$ws.Range("F27").Value = "https://www.iso.org/obp/ui/#search`n"

# Row 28
$ws.Range("A28").Value = "vars:GrowthStageName"
$ws.Range("B28").Value = "GrowthStageName"
$ws.Range("F28").Value = "https://en.wikipedia.org/wiki/BBCH-scale_(potato)"

# Row 29
$ws.Range("A29").Value = "vars:SeverityCategoryName"
$ws.Range("B29").Value = "SeverityCategoryName"

# New rows 93 and 94 (template "vars:" rows, matching existing blank rows like row 92)
$ws.Range("A93").Value = "vars:"
$ws.Range("B93:T93").Value = ""
$ws.Range("A94").Value = "vars:"
$ws.Range("B94:T94").Value = ""
